# Add via upload: refresh season-header labels with a trailing "年" (year)
# character, and move the active cell selection to I1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Boston Celtics")

$ws.Range("B1").Value = "21-22年"
$ws.Range("C1").Value = "20-21年"
$ws.Range("D1").Value = "19-20年"
$ws.Range("E1").Value = "18-19年"
$ws.Range("F1").Value = "17-18年"
$ws.Range("G1").Value = "16-17年"
$ws.Range("H1").Value = "15-16年"
$ws.Range("I1").Value = "14-15年"

$ws.Range("I1").Select()
